# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh
# the handoff timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: B2 (zh-cn status), C2 (de-de status), D2 (latest handoff date)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-27-19 12:27:23"

# zh-cn sheet: C2 (status), E2 (latest handoff datetime)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-19 12:27:20"

# de-de sheet: C2 (status), E2 (latest handoff datetime)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-19 12:27:23"
